# Add two new columns "I0" (column I) and "IF" (column J) to the sheet,
# mirroring the header style used by the existing header row (e.g. H1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Reuse the same formatting as the existing header cells (bold, bordered,
# centered) so the new header cells share the same style index.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for rows 2-64 (63 data rows)
$iVals = @(9,8,9,9,9,9,9,9,10,9,9,9,9,8,9,9,9,8,9,9,9,9,9,9,9,9,8,9,9,9,8,9,9,9,8,9,9,9,10,9,8,8,8,9,9,10,8,9,9,9,9,9,8,9,9,9,6,5,8,7,7,6,8)
$jVals = @(9,9,9,9,9,9,9,9,10,10,10,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,8,9,9,9,9,9,9,9,9,10,9,9,11,9,8,9,9,9,9,10,9,9,9,9,9,9,9,10,9,9,6,5,8,7,7,6,8)

for ($r = 0; $r -lt 63; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$r]
    $ws.Cells.Item($row, 10).Value = $jVals[$r]
}
